# Apply the "two-digit number multiplied by two-digit number" worksheet update.
# Each pair below is (old text, new text); every old text is unique in the
# document, so a simple Find/Replace (ReplaceAll) for each is sufficient.

$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-30 Wednesday", "2025-07-31 Thursday"),
    @("47×91=", "17×11="),
    @("33×11=", "34×13="),
    @("67×69=", "77×76="),
    @("55×15=", "93×21="),
    @("21×78=", "89×54="),
    @("56×71=", "69×26="),
    @("73×96=", "29×25="),
    @("82×14=", "28×44="),
    @("27×99=", "98×97="),
    @("23×12=", "92×65="),
    @("84×73=", "45×40="),
    @("15×38=", "93×89="),
    @("20×42=", "67×46="),
    @("67×51=", "45×89="),
    @("91×74=", "65×41="),
    @("14×48=", "51×84="),
    @("27×26=", "73×65="),
    @("87×32=", "81×99="),
    @("42×32=", "60×94="),
    @("36×39=", "69×51="),
    @("17×49=", "34×70="),
    @("87×70=", "45×61="),
    @("60×85=", "74×13="),
    @("38×86=", "90×65="),
    @("71×92=", "32×59=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
